# Monthly rollover update
#   - "VENTAS POR GRUPO": zero-out the per-category current-month breakdown
#     (columns C..R, rows 2..55) and refresh the "N de 54" summary row (56).
#   - "VENTA MENSUAL": shift the 4 rolling month columns (C..F) one month to
#     the left (abril/mayo/junio/julio -> mayo/junio/julio/agosto), dropping
#     the oldest month and leaving the new rightmost month blank (0), then
#     recompute the column totals in row 56.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Narrow the LED column (J) a bit, matching the refreshed content width.
$wsGrupo.Columns.Item(10).ColumnWidth = 9

$firstRow = 2
$lastRow = 55
$firstCol = 3   # C
$lastCol = 18   # R

# Zero out every numeric data cell in the category breakdown.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $wsGrupo.Cells.Item($r, $c).Value2 = 0
    }
}

# Refresh the "x de 54" tally row: every column now has zero non-zero rows.
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $wsGrupo.Cells.Item(56, $c).Value2 = "0 de 54"
}

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths follow the refreshed header/content widths.
$wsMensual.Columns.Item(3).ColumnWidth = 14   # C
$wsMensual.Columns.Item(6).ColumnWidth = 12   # F

# Shift the month headers one month forward.
$wsMensual.Cells.Item(1, 3).Value2 = "mayo"
$wsMensual.Cells.Item(1, 4).Value2 = "junio"
$wsMensual.Cells.Item(1, 5).Value2 = "julio"
$wsMensual.Cells.Item(1, 6).Value2 = "agosto"

$dataFirstRow = 2
$dataLastRow = 55

for ($r = $dataFirstRow; $r -le $dataLastRow; $r++) {
    $oldD = $wsMensual.Cells.Item($r, 4).Value2
    $oldE = $wsMensual.Cells.Item($r, 5).Value2
    $oldF = $wsMensual.Cells.Item($r, 6).Value2

    $wsMensual.Cells.Item($r, 3).Value2 = $oldD   # C <- D
    $wsMensual.Cells.Item($r, 4).Value2 = $oldE   # D <- E
    $wsMensual.Cells.Item($r, 5).Value2 = $oldF   # E <- F
    $wsMensual.Cells.Item($r, 6).Value2 = 0       # F <- blank new month
}

# AGUIMPORT-AGUILAR IMPORTACIONES S.A.S. (row 5) got a corrected June figure
# from the upstream source feed (not a pure carry-over of the old July cell).
$wsMensual.Cells.Item(5, 5).Value2 = 11268.16

# Recompute the column totals (row 56) for the shifted data. Use a
# Neumaier (compensated) running sum so the totals line up with the
# upstream report's precision instead of drifting in the last bit.
for ($c = 3; $c -le 6; $c++) {
    $s = 0.0
    $comp = 0.0
    for ($r = $dataFirstRow; $r -le $dataLastRow; $r++) {
        $v = $wsMensual.Cells.Item($r, $c).Value2
        $t = $s + $v
        if ([Math]::Abs($s) -ge [Math]::Abs($v)) {
            $comp = $comp + (($s - $t) + $v)
        } else {
            $comp = $comp + (($v - $t) + $s)
        }
        $s = $t
    }
    $wsMensual.Cells.Item(56, $c).Value2 = $s + $comp
}
